$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "modified the move message's body" -- the RelivePos (move/respawn position)
# column values for rows 11-17 are unified to a single coordinate "55,110,0".
$ws.Range("G11:G17").Value = "55,110,0"

# Selection moves from the old multi-cell range G11:G17 to a single cell G15.
$ws.Range("G15").Select()
